$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("SoFCtMbCtPR ")

$ws.Range("B2:B24").Value = 1.075
